$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "target" column (G) values were renamed from the full word "deuteron"
# to the shorthand "d" for every data row.
for ($r = 2; $r -le 13; $r++) {
    $ws.Range("G$r").Value = "d"
}

# The header row (row 1) was made bold.
$ws.Range("A1:K1").Font.Bold = $true

# The active selection when the file was last saved moved to E20.
$ws.Range("E20").Select() | Out-Null
